# Update the cryptos price (column D) and 1h volume (column E) figures to
# match the latest scrape, per the "Updated cryptos list ... with GitHub
# Actions" commit. Every D/E cell in this sheet is stored as text, so for
# any new D value that looks like a plain number we force the cell's
# number format to Text ("@") first -- otherwise Excel's COM layer would
# silently reinterpret e.g. "303.59" as the number 303.59 and drop the
# original formatting/precision semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "23.231.49";    E = "  +0.90%  " },
    @{ Row = 3;  D = "1.603.00";     E = "  +0.01%  " },
    @{ Row = 4;  D = $null;          E = "  -0.04%  " },
    @{ Row = 5;  D = $null;          E = "  -0.11%  " },
    @{ Row = 6;  D = "303.59";       E = "  +0.79%  " },
    @{ Row = 7;  D = "0.3766";       E = "  -0.32%  " },
    @{ Row = 8;  D = "51.86";        E = "  +4.66%  " },
    @{ Row = 9;  D = "0.3632";       E = "  +0.20%  " },
    @{ Row = 10; D = $null;          E = "  +1.24%  " },
    @{ Row = 11; D = $null;          E = "  -0.03%  " },
    @{ Row = 12; D = "0.08123";      E = "  +0.12%  " },
    @{ Row = 13; D = "22.80";        E = "  +0.03%  " },
    @{ Row = 14; D = "6.594";        E = "  +0.14%  " },
    @{ Row = 15; D = "7.417";        E = "  +0.18%  " },
    @{ Row = 16; D = "0.00001248";   E = "  +0.32%  " },
    @{ Row = 17; D = "1.605.47";     E = "  +0.43%  " },
    @{ Row = 18; D = "94.07";        E = "  +2.15%  " },
    @{ Row = 19; D = "0.06920";      E = "  +0.60%  " },
    @{ Row = 20; D = "18.15";        E = "  -0.30%  " },
    @{ Row = 21; D = "6.527";        E = "  -0.49%  " },
    @{ Row = 22; D = $null;          E = "  -0.16%  " },
    @{ Row = 23; D = $null;          E = "  -1.60%  " },
    @{ Row = 24; D = "23.225.14";    E = "  +0.85%  " },
    @{ Row = 25; D = $null;          E = "  +8.48%  " },
    @{ Row = 26; D = "2.377";        E = "  +0.90%  " },
    @{ Row = 27; D = "21.21";        E = "  +0.66%  " },
    @{ Row = 28; D = "150.13";       E = "  -0.19%  " },
    @{ Row = 29; D = "5.259";        E = "  -0.03%  " },
    @{ Row = 30; D = "134.67";       E = "  +0.87%  " },
    @{ Row = 31; D = "2.400";        E = "  +4.20%  " },
    @{ Row = 32; D = "6.739";        E = "  -1.02%  " },
    @{ Row = 33; D = "1.782.86";     E = "  -0.01%  " },
    @{ Row = 34; D = "0.9601";       E = "  -0.43%  " },
    @{ Row = 35; D = "0.07496";      E = "  -1.69%  " },
    @{ Row = 36; D = "0.02745";      E = "  +1.41%  " },
    @{ Row = 37; D = "10.29";        E = "  -0.86%  " },
    @{ Row = 38; D = "0.2535";       E = "  +0.06%  " },
    @{ Row = 39; D = "6.120";        E = "  -2.56%  " },
    @{ Row = 40; D = "0.08822";      E = "  -0.22%  " },
    @{ Row = 41; D = "1.393";        E = "  +2.08%  " },
    @{ Row = 42; D = "0.7108";       E = "  +0.79%  " },
    @{ Row = 43; D = "12.50";        E = "  -0.10%  " },
    @{ Row = 44; D = "15.66";        E = "  +3.15%  " },
    @{ Row = 45; D = "0.6545";       E = "  -0.94%  " },
    @{ Row = 46; D = "2.315";        E = "  +0.10%  " },
    @{ Row = 47; D = "0.9996";       E = "  -0.10%  " },
    @{ Row = 48; D = $null;          E = "  +0.67%  " },
    @{ Row = 49; D = "132.64";       E = "  +0.09%  " },
    @{ Row = 50; D = "0.07948";      E = "  +0.47%  " },
    @{ Row = 51; D = $null;          E = "  -1.56%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        # Numeric-looking strings (e.g. "303.59") must stay text, just like
        # every other cell in this column -- pin the format to Text first,
        # otherwise Excel auto-converts the literal to a real number.
        if ($u.D -match '^-?[0-9]+(\.[0-9]+)?$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
